$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column P (appliedMarkCode) first so column letters for G stay valid
$ws.Range("P1").EntireColumn.Delete()

# Delete column G (sourceOfFishSite)
$ws.Range("G1").EntireColumn.Delete()
